# "add sign VTHO transactions console"
# Update the three VTHO transaction rows: bump the amount from 1.88 to
# 188000 and point them at the new block reference, then re-select the
# console's next input cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - also refresh the "To" address text (same value, rewritten)
$ws.Range("A3").Value = "0xD3EF28DF6b553eD2fc47259E8134319cB1121A2A"
$ws.Range("B3").Value = 188000
$ws.Range("D3").Value = "0x00015e41be43bb95"

# Row 4
$ws.Range("B4").Value = 188000
$ws.Range("D4").Value = "0x00015e41be43bb95"

# Row 5
$ws.Range("B5").Value = 188000
$ws.Range("D5").Value = "0x00015e41be43bb95"

# Move the active selection to the next console input cell
$ws.Range("D12").Select()
